# Update stats for 2025-11 (row 24 in Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B24").Value = 6348
$ws.Range("D24").Value = 5948908
$ws.Range("E24").Value = 937.1310649023314
$ws.Range("F24").Value = 8.216842823048065
$ws.Range("H24").Value = 26.02091441460907
